# Update cryptocurrency price/volume data (auto-refresh commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.492.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.566.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  -1.32%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +4.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.05"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("E10").Value = "  -1.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0593"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.789.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.584.50"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.491.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.53%  "
$ws.Range("E18").Value = "  -1.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "228.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.86%  "
$ws.Range("E21").Value = "  -2.63%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").Value = "  -6.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "150.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.81%  "
$ws.Range("E29").Value = "  -4.21%  "
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0467"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.90%  "
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.392.69"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.28%  "
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0166"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.537"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.78%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.788"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.48%  "
$ws.Range("E45").Value = "  -4.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.975"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "62.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.701.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "85.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.04%  "
$ws.Range("E50").Value = "  -3.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0525"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.51%  "
